# Pedidos.xlsx update
# ---------------------------------------------------------------------
# - drop two shipments that are no longer pending (80265160/30018-KMT-I
#   and the first 80265942/20637-TDK-I line)
# - change the Quantidade column to a plain integer look instead of the
#   old "#,##0.000" one
# - add the newly arrived shipments (remessas 80266501 .. 80266508)
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents for rows 2-20 (Remessa / Material / Quantidade).
$remessas = @(
    "80265942","80265942","80265942","80265942",
    "80265944","80266080","80266081","80266324","80266481","80266490",
    "80266491","80266492","80266501","80266502","80266504","80266505",
    "80266506","80266507","80266508"
)
$materiais = @(
    "20850-FUZ-I","20869-FUZ-I","20853-FUZ-I","60234-WUE-I",
    "20389-DCC-I","60233-STM-I","40069-TDK-I","40193-TDK-N","60208-STM-I","30253-OSR-I",
    "10255-ARI-I","10000-LDG-I","10145-ARI-I","10000-LDG-I","12003-KRO-I","10636-ARI-I",
    "20935-CTY-I","10000-LDG-I","10000-LDG-I"
)
$quantidades = @(600,36000,1900,1500,5000,2000,1064,15000,3,384000,1,2,1,1,4,1,2,1,1)

$firstRow = 2
$lastRow  = $firstRow + $remessas.Length - 1

# Columns A & B hold values that look like numbers ("80265942", "10000-LDG-I"
# style codes, etc.) but must stay text, exactly like they already are
# elsewhere on the sheet. Force a Text format before typing them in so Excel
# doesn't silently convert them to numbers.
$colA = $ws.Range($ws.Cells.Item($firstRow,1), $ws.Cells.Item($lastRow,1))
$colB = $ws.Range($ws.Cells.Item($firstRow,2), $ws.Cells.Item($lastRow,2))
$colA.NumberFormat = "@"
$colB.NumberFormat = "@"

for ($i = 0; $i -lt $remessas.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 1).Value = $remessas[$i]
    $ws.Cells.Item($row, 2).Value = $materiais[$i]
    $ws.Cells.Item($row, 3).Value = $quantidades[$i]
}

# Re-apply the normal (General) look used by the rest of column A & B, and
# the plain-integer look already used further down column C, by copying the
# formatting from cells that already carry it correctly.
$ws.Range("A21").Copy()
$ws.Range($ws.Cells.Item($firstRow,1), $ws.Cells.Item($lastRow,1)).PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range($ws.Cells.Item($firstRow,2), $ws.Cells.Item($lastRow,2)).PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range($ws.Cells.Item($firstRow,3), $ws.Cells.Item($lastRow,3)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Cosmetic: cursor/selection left on E10 by the author.
$ws.Range("E10").Select()

Write-Host "done"
